$d = $word.ActiveDocument

# Portuguese paragraph: split after "Métodos de difração de raios X. " (keep trailing space on first part)
$d.Content.Find.Execute(
    "Métodos de difração de raios X. Preparação materialográfica",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Métodos de difração de raios X. ^lPreparação materialográfica", 2
)

# English paragraph: split right after "Methods of X-ray diffraction." (no space)
$d.Content.Find.Execute(
    "Methods of X-ray diffraction.Materialographic sample preparation",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Methods of X-ray diffraction.^lMaterialographic sample preparation", 2
)
